$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Run" column values for rows 3-5 from "Yes" to "No"
$ws.Range("F3").Value = "No"
$ws.Range("F4").Value = "No"
$ws.Range("F5").Value = "No"

# Move selection to match the new active cell/range (F4:F5, active cell F4)
$ws.Range("F4:F5").Select()
